# Insert a new data row at row 527 (pushing the existing rows 527-574 down to
# 528-575) and populate it with the new "Tercera" quality record dated
# 2023-09-25. All other columns reuse the constant metadata shared by every
# row in this dataset (Mercado, Region, etc.).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(527).Insert()

$ws.Cells.Item(527, 1).Value  = 3
$ws.Cells.Item(527, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(527, 3).Value  = "Coquimbo"
$ws.Cells.Item(527, 4).Value  = "2023-09-25"
$ws.Cells.Item(527, 5).Value  = 5
$ws.Cells.Item(527, 6).Value  = "Fruta"
$ws.Cells.Item(527, 7).Value  = 100101
$ws.Cells.Item(527, 8).Value  = "Berries"
$ws.Cells.Item(527, 9).Value  = 100112025
$ws.Cells.Item(527, 10).Value = "Frutilla"
$ws.Cells.Item(527, 11).Value = "Sin especificar"
$ws.Cells.Item(527, 12).Value = "Tercera"
$ws.Cells.Item(527, 13).Value = 54
$ws.Cells.Item(527, 14).Value = 12000
$ws.Cells.Item(527, 15).Value = 12000
$ws.Cells.Item(527, 16).Value = 12000
$ws.Cells.Item(527, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(527, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(527, 19).Value = 1714
$ws.Cells.Item(527, 20).Value = 7
